$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.45200502670175
$ws.Range("C2").Value = 9.545637825307402
$ws.Range("D2").Value = 5.961665458861129
$ws.Range("E2").Value = 11.65020442559076
$ws.Range("G2").Value = 20.77794550640796
$ws.Range("H2").Value = 12.12517691608559
$ws.Range("I2").Value = 16.68176608603568
$ws.Range("M2").Value = 14.43847468672994
$ws.Range("O2").Value = 17.40208327527908

$ws.Range("B3").Value = 11.75953280826774
$ws.Range("C3").Value = 9.059483628876681
$ws.Range("D3").Value = 5.838380972306253
$ws.Range("E3").Value = 11.59103848945659
$ws.Range("G3").Value = 20.75410458154744
$ws.Range("H3").Value = 12.17900810912763
$ws.Range("I3").Value = 16.81902821451002
$ws.Range("M3").Value = 14.10383819417179
$ws.Range("O3").Value = 17.47408869764488

$ws.Range("B4").Value = 11.31237687855554
$ws.Range("C4").Value = 8.745855238147488
$ws.Range("D4").Value = 5.763126903305853
$ws.Range("E4").Value = 11.56003877484673
$ws.Range("G4").Value = 20.7516498165664
$ws.Range("H4").Value = 12.21493039547988
$ws.Range("I4").Value = 16.90826100312615
$ws.Range("M4").Value = 13.89682104021401
$ws.Range("O4").Value = 17.52433188094259

$ws.Range("B5").Value = 11.12476315751048
$ws.Range("C5").Value = 8.614335992735315
$ws.Range("D5").Value = 5.732620069700396
$ws.Range("E5").Value = 11.5487539422179
$ws.Range("G5").Value = 20.75370249919656
$ws.Range("H5").Value = 12.23028889790592
$ws.Range("I5").Value = 16.94586876144686
$ws.Range("M5").Value = 13.81219183620578
$ws.Range("O5").Value = 17.54631438778335

$ws.Range("B6").Value = 11.09328866222862
$ws.Range("C6").Value = 8.592276060601147
$ws.Range("D6").Value = 5.72756551200774
$ws.Range("E6").Value = 11.54696168430116
$ws.Range("G6").Value = 20.75422732649898
$ws.Range("H6").Value = 12.23288259366086
$ws.Range("I6").Value = 16.95218865309666
$ws.Range("M6").Value = 13.79812666908578
$ws.Range("O6").Value = 17.55005538574761

$ws.Range("B7").Value = 11.30986829996507
$ws.Range("C7").Value = 8.744096425041308
$ws.Range("D7").Value = 5.762714765983187
$ws.Range("E7").Value = 11.55988111809642
$ws.Range("G7").Value = 20.75166515621481
$ws.Range("H7").Value = 12.21513461316498
$ws.Range("I7").Value = 16.9087631559654
$ws.Range("M7").Value = 13.89568062267203
$ws.Range("O7").Value = 17.52462225062184

$ws.Range("B8").Value = 12.2178723399714
$ws.Range("C8").Value = 9.381200923834058
$ws.Range("D8").Value = 5.919093796701369
$ws.Range("E8").Value = 11.62870372559609
$ws.Range("G8").Value = 20.76719164654418
$ws.Range("H8").Value = 12.14314126237665
$ws.Range("I8").Value = 16.72806546450873
$ws.Range("M8").Value = 14.32348496103234
$ws.Range("O8").Value = 17.42565377121143

$ws.Range("B9").Value = 13.81958544903858
$ws.Range("C9").Value = 10.50734608585027
$ws.Range("D9").Value = 6.227188997711079
$ws.Range("E9").Value = 11.80545701758506
$ws.Range("G9").Value = 20.89452249513073
$ws.Range("H9").Value = 12.0248096041831
$ws.Range("I9").Value = 16.41306105520731
$ws.Range("M9").Value = 15.14503066546561
$ws.Range("O9").Value = 17.27980750729611

$ws.Range("B10").Value = 14.88271011773838
$ws.Range("C10").Value = 11.2562887167191
$ws.Range("D10").Value = 6.45173501302622
$ws.Range("E10").Value = 11.95995374227366
$ws.Range("G10").Value = 21.04703960877555
$ws.Range("H10").Value = 11.95190456480763
$ws.Range("I10").Value = 16.20567217303338
$ws.Range("M10").Value = 15.73161841974233
$ws.Range("O10").Value = 17.20254272075706

$ws.Range("B11").Value = 15.34102974132382
$ws.Range("C11").Value = 11.57948718242632
$ws.Range("D11").Value = 6.552961231042524
$ws.Range("E11").Value = 12.03535502324321
$ws.Range("G11").Value = 21.1290900965726
$ws.Range("H11").Value = 11.92180776706924
$ws.Range("I11").Value = 16.11656277979843
$ws.Range("M11").Value = 15.99355211437488
$ws.Range("O11").Value = 17.17397833176275

$ws.Range("B12").Value = 15.51090712242015
$ws.Range("C12").Value = 11.69932809685092
$ws.Range("D12").Value = 6.591119241447887
$ws.Range("E12").Value = 12.06462083622793
$ws.Range("G12").Value = 21.16196252437699
$ws.Range("H12").Value = 11.91085394672758
$ws.Range("I12").Value = 16.08357367090982
$ws.Range("M12").Value = 16.09194067770742
$ws.Range("O12").Value = 17.16411558129834

$ws.Range("B13").Value = 15.4744851165503
$ws.Range("C13").Value = 11.67363192506541
$ws.Range("D13").Value = 6.582909615135469
$ws.Range("E13").Value = 12.05828660122183
$ws.Range("G13").Value = 21.15480310613198
$ws.Range("H13").Value = 11.91319330554951
$ws.Range("I13").Value = 16.090644863847
$ws.Range("M13").Value = 16.07078794846697
$ws.Range("O13").Value = 17.1661971637304

$ws.Range("B14").Value = 15.35507956901145
$ws.Range("C14").Value = 11.58939777174285
$ws.Range("D14").Value = 6.556104234207076
$ws.Range("E14").Value = 12.03774859038798
$ws.Range("G14").Value = 21.131758556979
$ws.Range("H14").Value = 11.92089769622295
$ws.Range("I14").Value = 16.11383360765866
$ws.Range("M14").Value = 16.00166314968356
$ws.Range("O14").Value = 17.17314775893978

$ws.Range("B15").Value = 15.28146017893926
$ws.Range("C15").Value = 11.53746938245318
$ws.Range("D15").Value = 6.539661287080397
$ws.Range("E15").Value = 12.02526059257817
$ws.Range("G15").Value = 21.11787704829203
$ws.Range("H15").Value = 11.92567463263043
$ws.Range("I15").Value = 16.12813573364011
$ws.Range("M15").Value = 15.9592153476878
$ws.Range("O15").Value = 17.17752963910448

$ws.Range("B16").Value = 14.85224507496684
$ws.Range("C16").Value = 11.23481199915265
$ws.Range("D16").Value = 6.445097486719845
$ws.Range("E16").Value = 11.95512717640861
$ws.Range("G16").Value = 21.04193075834505
$ws.Range("H16").Value = 11.9539333699766
$ws.Range("I16").Value = 16.2116011642288
$ws.Range("M16").Value = 15.71439345525602
$ws.Range("O16").Value = 17.20454258229793

$ws.Range("B17").Value = 14.58242351761402
$ws.Range("C17").Value = 11.04463520326842
$ws.Range("D17").Value = 6.386819321025457
$ws.Range("E17").Value = 11.91339805104934
$ws.Range("G17").Value = 20.99857251275006
$ws.Range("H17").Value = 11.97205642365266
$ws.Range("I17").Value = 16.26414627320821
$ws.Range("M17").Value = 15.56287678235235
$ws.Range("O17").Value = 17.22280570967702

$ws.Range("B18").Value = 14.42484951180465
$ws.Range("C18").Value = 10.93360478031482
$ws.Range("D18").Value = 6.353214927589508
$ws.Range("E18").Value = 11.88988000616639
$ws.Range("G18").Value = 20.97482808312922
$ws.Range("H18").Value = 11.98276898738982
$ws.Range("I18").Value = 16.29486122733162
$ws.Range("M18").Value = 15.47527162536438
$ws.Range("O18").Value = 17.23392947035048

$ws.Range("B19").Value = 14.37109031857941
$ws.Range("C19").Value = 10.89573018582882
$ws.Range("D19").Value = 6.341823884983209
$ws.Range("E19").Value = 11.88200090705901
$ws.Range("G19").Value = 20.96699424049186
$ws.Range("H19").Value = 11.98644559434544
$ws.Range("I19").Value = 16.30534528852403
$ws.Range("M19").Value = 15.44553456097942
$ws.Range("O19").Value = 17.23780192012025

$ws.Range("B20").Value = 14.61139310332031
$ws.Range("C20").Value = 11.06505042418351
$ws.Range("D20").Value = 6.393032169700414
$ws.Range("E20").Value = 11.91779031498226
$ws.Range("G20").Value = 21.00306461788384
$ws.Range("H20").Value = 11.97009730360634
$ws.Range("I20").Value = 16.25850177995846
$ws.Range("M20").Value = 15.57905397051218
$ws.Range("O20").Value = 17.22079742007726

$ws.Range("B21").Value = 15.39025197149862
$ws.Range("C21").Value = 11.61420873084959
$ws.Range("D21").Value = 6.563982671349586
$ws.Range("E21").Value = 12.0437619471506
$ws.Range("G21").Value = 21.13847858957386
$ws.Range("H21").Value = 11.91862268619288
$ws.Range("I21").Value = 16.10700201107176
$ws.Range("M21").Value = 16.02198920800327
$ws.Range("O21").Value = 17.17108025531162

$ws.Range("B22").Value = 15.87782925075119
$ws.Range("C22").Value = 11.95825790586714
$ws.Range("D22").Value = 6.674675976744974
$ws.Range("E22").Value = 12.13023535254203
$ws.Range("G22").Value = 21.23746801238928
$ws.Range("H22").Value = 11.8875651777201
$ws.Range("I22").Value = 16.01238818931243
$ws.Range("M22").Value = 16.30677096815672
$ws.Range("O22").Value = 17.1441506241836

$ws.Range("B23").Value = 15.61957370546807
$ws.Range("C23").Value = 11.77600034911168
$ws.Range("D23").Value = 6.615704512758569
$ws.Range("E23").Value = 12.08371191506409
$ws.Range("G23").Value = 21.18368381968487
$ws.Range("H23").Value = 11.90390405036578
$ws.Range("I23").Value = 16.06248198567041
$ws.Range("M23").Value = 16.15523714509209
$ws.Range("O23").Value = 17.15801217060283

$ws.Range("B24").Value = 14.59830358170995
$ws.Range("C24").Value = 11.0558259807798
$ws.Range("D24").Value = 6.390223647550605
$ws.Range("E24").Value = 11.9158030986554
$ws.Range("G24").Value = 21.00103005131505
$ws.Range("H24").Value = 11.97098210804646
$ws.Range("I24").Value = 16.26105207909185
$ws.Range("M24").Value = 15.57174180206741
$ws.Range("O24").Value = 17.22170342446033

$ws.Range("B25").Value = 13.40597373161038
$ws.Range("C25").Value = 10.2162730341465
$ws.Range("D25").Value = 6.143971464634494
$ws.Range("E25").Value = 11.75324296715695
$ws.Range("G25").Value = 20.8496920468704
$ws.Range("H25").Value = 12.05436463935812
$ws.Range("I25").Value = 16.4940605171451
$ws.Range("M25").Value = 14.92532700853111
$ws.Range("O25").Value = 17.31404760726798
